$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7996079921722412
$ws.Range("B1").Value = 1.715750694274902
$ws.Range("C1").Value = 4.790500640869141
$ws.Range("D1").Value = 5.086749076843262
$ws.Range("E1").Value = 1.767454147338867
